$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new API row (row 9): Method / URL / description
$ws.Cells.Item(9, 2).Value = "POST"
$ws.Cells.Item(9, 3).Value = "/api3/create_diagnosis/"
$ws.Cells.Item(9, 4).Value = "진료 데이터 추가"

# Update the active selection to C8
[void]$ws.Range("C8").Select()
